$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for the affected rows
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -3
$ws.Range("F9").Value = 4
$ws.Range("F22").Value = 0
